$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1127464861192516
$ws.Range("B3").Value = 0.07426135796710345
$ws.Range("H3").Value = 0.1870078440863551
$ws.Range("B4").Value = 0.07989740725516165
$ws.Range("H4").Value = 0.1926438933744133
$ws.Range("B5").Value = 0.05793317263531943
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 0.1706796587545711
$ws.Range("B6").Value = 0.04473943535494208
$ws.Range("C6").Value = 0.005325043559503317
$ws.Range("D6").Value = 3.648885329716293
$ws.Range("E6").Value = 0.01345280638785308
$ws.Range("F6").Value = 0.03425518317390631
$ws.Range("G6").Value = 0.05522368753597827
$ws.Range("H6").Value = 0.1574859214741937
$ws.Range("B7").Value = 0.02903735774562823
$ws.Range("C7").Value = 0.003044100755573473
$ws.Range("D7").Value = 2.139757338538166
$ws.Range("E7").Value = 0.00005900886818307892
$ws.Range("F7").Value = 0.02306212863388198
$ws.Range("G7").Value = 0.03501258685737458
$ws.Range("H7").Value = 0.1417838438648799
$ws.Range("B8").Value = 0.02866960613684441
$ws.Range("C8").Value = 0.002732353517780752
$ws.Range("D8").Value = 2.210747649709623
$ws.Range("E8").Value = 0.0000003076690044233385
$ws.Range("F8").Value = 0.02330582168138613
$ws.Range("G8").Value = 0.03403339059230218
$ws.Range("H8").Value = 0.141416092256096
$ws.Range("B9").Value = 0.0138662000905528
$ws.Range("C9").Value = 0.003605075256054655
$ws.Range("D9").Value = 0.005206977931849108
$ws.Range("E9").Value = 0.00001778910002200519
$ws.Range("F9").Value = 0.006790550059476552
$ws.Range("G9").Value = 0.02094185012162898
$ws.Range("H9").Value = 0.1266126862098044
$ws.Range("B10").Value = 0.00001859951780417363
$ws.Range("C10").Value = 0.001127331343166753
$ws.Range("D10").Value = -1.00913892949151
$ws.Range("E10").Value = 0.0000001995100196980731
$ws.Range("F10").Value = -0.002198755385376016
$ws.Range("G10").Value = 0.002235954420984398
$ws.Range("H10").Value = 0.1127650856370558
$ws.Range("B11").Value = 0.02657163319265903
$ws.Range("H11").Value = 0.1393181193119107
$ws.Range("B12").Value = 0.04287082235326538
$ws.Range("H12").Value = 0.155617308472517
$ws.Range("B13").Value = 0.0491633941321979
$ws.Range("H13").Value = 0.1619098802514495
$ws.Range("B14").Value = 0.0514568738883066
$ws.Range("C14").Value = 0.008987261085764024
$ws.Range("D14").Value = 9.198863409012212
$ws.Range("E14").Value = 0.0444454659280662
$ws.Range("F14").Value = 0.03379007168371145
$ws.Range("G14").Value = 0.06912367609290182
$ws.Range("H14").Value = 0.1642033600075582
$ws.Range("B15").Value = 0.05828135157323943
$ws.Range("H15").Value = 0.1710278376924911
$ws.Range("B16").Value = 0.0619259653860847
$ws.Range("C16").Value = 0.009099232169588472
$ws.Range("D16").Value = 9.907563843853747
$ws.Range("E16").Value = 0.04583561235697588
$ws.Range("F16").Value = 0.04404726587570973
$ws.Range("G16").Value = 0.07980466489645974
$ws.Range("H16").Value = 0.1746724515053363
$ws.Range("B17").Value = 0.06708972670187829
$ws.Range("C17").Value = 0.009059498303894638
$ws.Range("D17").Value = 10.28598447811154
$ws.Range("E17").Value = 0.03945756713915682
$ws.Range("F17").Value = 0.0492684602889722
$ws.Range("G17").Value = 0.08491099311478427
$ws.Range("H17").Value = 0.1798362128211299
$ws.Range("B18").Value = -0.1127464861192516
$ws.Range("C18").Value = 0.01268551493633767
$ws.Range("D18").Value = -16.06863000417736
$ws.Range("E18").Value = 0.01260794727490581
$ws.Range("F18").Value = -0.1376583015129513
$ws.Range("G18").Value = -0.08783467072555183
$ws.Range("B19").Value = 0.06296001873867402
$ws.Range("C19").Value = 0.008772607598267574
$ws.Range("D19").Value = 10.12573975635952
$ws.Range("E19").Value = 0.04611348190673749
$ws.Range("F19").Value = 0.04569548489236878
$ws.Range("G19").Value = 0.08022455258497924
$ws.Range("H19").Value = 0.1757065048579257
$ws.Range("B20").Value = 0.06575380948743922
$ws.Range("C20").Value = 0.008779186937611547
$ws.Range("D20").Value = 10.60119039919291
$ws.Range("E20").Value = 0.0454350173039485
$ws.Range("F20").Value = 0.04853220242406418
$ws.Range("G20").Value = 0.08297541655081425
$ws.Range("H20").Value = 0.1785002956066908
$ws.Range("B21").Value = 0.06845930811803592
$ws.Range("C21").Value = 0.009017170268084507
$ws.Range("D21").Value = 11.11654896095957
$ws.Range("E21").Value = 0.04765362754886252
$ws.Range("F21").Value = 0.0507228984519586
$ws.Range("G21").Value = 0.08619571778411324
$ws.Range("H21").Value = 0.1812057942372876
$ws.Range("B22").Value = 0.06723363870094134
$ws.Range("C22").Value = 0.008570173494798092
$ws.Range("D22").Value = 10.64023457573836
$ws.Range("E22").Value = 0.04404701339746033
$ws.Range("F22").Value = 0.05042658461121194
$ws.Range("G22").Value = 0.08404069279067078
$ws.Range("H22").Value = 0.179980124820193
$ws.Range("B23").Value = 0.0709451441382204
$ws.Range("C23").Value = 0.009136867550119877
$ws.Range("D23").Value = 10.86079618772808
$ws.Range("E23").Value = 0.04601868827903455
$ws.Range("F23").Value = 0.05296322388597911
$ws.Range("G23").Value = 0.08892706439046179
$ws.Range("H23").Value = 0.183691630257472
$ws.Range("B24").Value = 0.06880958640594251
$ws.Range("C24").Value = 0.008419400757893004
$ws.Range("D24").Value = 10.97351591077867
$ws.Range("E24").Value = 0.04762091744852344
$ws.Range("F24").Value = 0.05228936199531305
$ws.Range("G24").Value = 0.08532981081657202
$ws.Range("H24").Value = 0.1815560725251941
$ws.Range("B25").Value = 0.07205723425062011
$ws.Range("C25").Value = 0.009293694449334101
$ws.Range("D25").Value = 11.22572366810779
$ws.Range("E25").Value = 0.05517921973289089
$ws.Range("F25").Value = 0.05373720482041149
$ws.Range("G25").Value = 0.09037726368082857
$ws.Range("H25").Value = 0.1848037203698717
$ws.Range("B26").Value = 0.07028253204454092
$ws.Range("C26").Value = 0.008480875819898617
$ws.Range("D26").Value = 11.25509295538604
$ws.Range("E26").Value = 0.06387091180657176
$ws.Range("F26").Value = 0.05365321403214315
$ws.Range("G26").Value = 0.08691185005693863
$ws.Range("H26").Value = 0.1830290181637926
$ws.Range("B27").Value = 0.07370065546307279
$ws.Range("C27").Value = 0.008935302678955991
$ws.Range("D27").Value = 10.84867078324668
$ws.Range("E27").Value = 0.06086329690069334
$ws.Range("F27").Value = 0.05615715177024259
$ws.Range("G27").Value = 0.09124415915590282
$ws.Range("H27").Value = 0.1864471415823244
$ws.Range("B28").Value = 0.07264075851552446
$ws.Range("C28").Value = 0.008723775816604214
$ws.Range("D28").Value = 10.75418622437872
$ws.Range("E28").Value = 0.08402750342740545
$ws.Range("F28").Value = 0.05551770485459395
$ws.Range("G28").Value = 0.08976381217645463
$ws.Range("H28").Value = 0.1853872446347761
$ws.Range("B29").Value = 0.003941038007569537
$ws.Range("C29").Value = 0.0014940422368962
$ws.Range("D29").Value = -0.7128851805879993
$ws.Range("E29").Value = 0.00008332458489466942
$ws.Range("F29").Value = 0.0009777423951667277
$ws.Range("G29").Value = 0.006904333619972322
$ws.Range("H29").Value = 0.1166875241268212
